# The post "「誰かを幸福にしなさい」" (row 498) was removed from the sheet.
# Deleting the entire row shifts every subsequent row up by one, which
# matches the target diff (old row 499 "「傷跡は星にかえなさい」" becomes
# the new row 498, ... , old row 685 "「昼寝の時間だ」" becomes the new
# row 684) and also updates the sheet's used-range dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(498).EntireRow.Delete()
